$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Sum total: x1:0.795|x2:0.957|x3:1.637|x4:1.035|x5:1.183|x6:0.966|x7:0.923|x8:1.041|x9:1.356|x10:0.970|x11:1.302|x12:1.121|x13:0.575|x14:0.574|x15:0.712|x16:0.854|x17:0.000"
$ws.Range("I2").Value = "0.941, 0.354"
$ws.Range("J2").Value = "68.280, 71.505"
$ws.Range("K2").Value = "82.343, 82.504"
